# Updates cryptos list values (Price and Volume(1h) columns) for rows 2-51
# matching the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.351.94"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "1.868.98"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'243.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'0.4699"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("D9").Value = "'0.06455"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").Value = "'22.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("D11").Value = "'0.07769"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "1.872.00"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "'95.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "'0.7206"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").Value = "'5.125"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "'279.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").Value = "30.338.58"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").Value = "'12.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.84%  "
$ws.Range("D19").Value = "'0.000007523"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "2.119.78"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'5.226"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").Value = "'6.230"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").Value = "'163.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").Value = "'9.034"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("D27").Value = "'18.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").Value = "'0.09592"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("E31").Value = "  -2.28%  "
$ws.Range("D32").Value = "'4.206"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("D34").Value = "'0.04811"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("D37").Value = "'2.711"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").Value = "'0.01876"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("D39").Value = "'2.811"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.04%  "
$ws.Range("D40").Value = "'6.211"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.04%  "
$ws.Range("D41").Value = "'74.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").Value = "'1.938"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.74%  "
$ws.Range("D43").Value = "'0.4217"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("D44").Value = "'0.9991"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "'0.8259"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("D46").Value = "'100.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("D47").Value = "'9.605"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.40%  "
$ws.Range("D48").Value = "'35.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").Value = "'6.943"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").Value = "'898.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("D51").Value = "'0.05715"
$ws.Range("D51").Style = "Normal"
